# Scen_NCAP_NAS.xlsx / NCAP_BND sheet
# Add a new "UP / CAP_BND / ELE_NEW_PV_GRND" bound row for year 2050 (value 18),
# inserted right after the existing 2040 row (row 33), pushing the
# "Ograniczenia dolne (wymuszenia)" block (and everything below it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33 - Excel shifts existing row 33.. down to 34.. and
# copies the formatting of the row above (row 32) onto the freshly inserted row,
# which is exactly the style pattern the new data row needs.
$ws.Rows("33:33").Insert()

# Match the row height used throughout this table (row-level height isn't
# carried over by Insert(), only the per-cell styles are).
$ws.Rows("33:33").RowHeight = 18.75

# Fill in the new row's values (continuing the ELE_NEW_PV_GRND / UP / CAP_BND series).
$ws.Range("B33").Value2 = "UP"
$ws.Range("C33").Value2 = "CAP_BND"
$ws.Range("D33").Value2 = 2050
$ws.Range("E33").Value2 = 18
$ws.Range("F33").Value2 = "ELE_NEW_PV_GRND"

# Leave the selection where the author left it when they saved the file.
$ws.Range("I24").Select() | Out-Null
